$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price values in column D (stored as text to preserve exact formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "266.62"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.34"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.120"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06101"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.574"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.491"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.358"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01346"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1585"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08086"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03384"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03205"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09217"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.742"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001632"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04653"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006462"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006143"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001068"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.726"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.262"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3276"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1243"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04599"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006992"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003899"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005815"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009899"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00001900"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01240"
